# Add ability to parse hyperlinks in cells.
# - New "URL" column (F) on the Posts sheet.
# - F2: a HYPERLINK() formula (the "friendly name" case).
# - F3: a real hyperlink inserted the way the Excel GUI does it
#       (Hyperlinks.Add against a cell that already holds the display text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Posts")

# Header for the new column.
$ws.Range("F1").Value = "URL"

# F2: hyperlink expressed as a formula - prioritizes the URL, but keeps
# the "friendly name" as the visible/cached text.
$ws.Range("F2").Formula = '=HYPERLINK("http://www.example.com/hyperlink-function", "This uses the HYPERLINK() function")'

# F3: hyperlink inserted the way the "Insert Hyperlink" GUI dialog does -
# put the friendly text in the cell first, then attach the hyperlink
# relationship to it (this is what applies the built-in "Hyperlink" cell
# style and records the relationship/hyperlink part).
$ws.Range("F3").Value = "This uses the hyperlink GUI option"

$hlStyle = $wb.Styles.Add("HyperlinkGui")
$hlStyle.Font.Underline = 2
$ws.Range("F3").Style = "HyperlinkGui"

$ws.Hyperlinks.Add($ws.Range("F3"), "http://www.example.com/hyperlink-gui")

# Match the updated selection left behind by the GUI hyperlink workflow.
$ws.Range("F3").Select()
